$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 155959
$ws.Range("C4").Value = 147053
$ws.Range("C5").Value = 8906
$ws.Range("C8").Value = 63.83
